$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Repayment schedule" sheet: widen column L (12) to match column I's
#    best-fit width (6.5703125 chars). The ColumnWidth COM property on this
#    host only resolves to a 1/6-character pixel grid, so 5.65 (which lands
#    on the nearest achievable grid value, 6.5) is used to get as close as
#    possible to the target width.
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns.Item(12).ColumnWidth = 5.65

# ---------------------------------------------------------------------------
# 2. "Transactions" sheet: update the selection and the transaction figures.
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

# Row 2
$wsTrans.Range("A2").Value = 61
$wsTrans.Range("J2").Value = 9133.2199999999993

# Row 3
$wsTrans.Range("A3").Value = 59
$wsTrans.Range("C3").Value = 42064
$wsTrans.Range("E3").Value = 963.77
$wsTrans.Range("F3").Value = 866.78
$wsTrans.Range("G3").Value = 96.99
$wsTrans.Range("J3").Value = 4133.22

# Row 4
$wsTrans.Range("A4").Value = 57

# Update the selection shown when the sheet is active (was C12, now A2:L4).
$wsTrans.Range("A2:L4").Select()
